$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.734.91'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '3.302.40'
$ws.Range("E3").Value = '  +6.03%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.10'
$ws.Range("E5").Value = '  +2.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.34'
$ws.Range("E6").Value = '  +4.81%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.301.31'
$ws.Range("E8").Value = '  +6.07%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +3.56%  '
$ws.Range("E11").Value = '  +4.91%  '
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.66'
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '3.844.93'
$ws.Range("E15").Value = '  +6.00%  '
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '3.302.37'
$ws.Range("E17").Value = '  +5.79%  '
$ws.Range("D18").Value = '63.806.71'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.68'
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").Value = '  +5.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.03'
$ws.Range("E23").Value = '  +5.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.51'
$ws.Range("E24").Value = '  +5.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.70'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +7.21%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +4.57%  '
$ws.Range("E31").Value = '  +5.07%  '
$ws.Range("E32").Value = '  +10.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.105'
$ws.Range("E33").Value = '  -1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.54'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.97'
$ws.Range("E36").Value = '  +4.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.10'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  +9.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0402'
$ws.Range("E39").Value = '  +4.61%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.060.88'
$ws.Range("E40").Value = '  +5.87%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '424.93'
$ws.Range("E41").Value = '  +2.22%  '
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  +3.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.112'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.19'
$ws.Range("E46").Value = '  +5.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.29'
$ws.Range("E47").Value = '  +4.17%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '124.24'
$ws.Range("E51").Value = '  +3.39%  '
